$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 0.9600122695046023),
    @(2, 0.9518562614462791),
    @(3, 0.9478176730952024),
    @(5, 0.9458168575092269),
    @(6, 0.9456571933548517),
    @(7, 0.9455623965112731),
    @(8, 0.9454888694856882),
    @(10, 0.9454018282780285)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
